{"js": "// Replace the date line and every two-digit-by-two-digit multiplication\n// problem in the practice sheet with the new values from the target\n// revision. Every \"before\" string below is unique within the document,\n// so a direct search()+insertText(\"Replace\") per pair is unambiguous.\nconst replacements = [\n  [\"2025-08-09 Saturday\", \"2025-08-10 Sunday\"],\n  [\"92\u00d740=\", \"61\u00d734=\"],\n  [\"14\u00d773=\", \"19\u00d745=\"],\n  [\"34\u00d716=\", \"78\u00d778=\"],\n  [\"49\u00d721=\", \"62\u00d713=\"],\n  [\"71\u00d768=\", \"92\u00d790=\"],\n  [\"39\u00d781=\", \"83\u00d772=\"],\n  [\"74\u00d782=\", \"95\u00d740=\"],\n  [\"85\u00d798=\", \"96\u00d721=\"],\n  [\"16\u00d766=\", \"40\u00d731=\"],\n  [\"34\u00d781=\", \"58\u00d756=\"],\n  [\"79\u00d773=\", \"42\u00d748=\"],\n  [\"62\u00d764=\", \"53\u00d732=\"],\n  [\"53\u00d737=\", \"81\u00d790=\"],\n  [\"47\u00d745=\", \"73\u00d785=\"],\n  [\"62\u00d721=\", \"18\u00d765=\"],\n  [\"45\u00d711=\", \"13\u00d771=\"],\n  [\"59\u00d756=\", \"49\u00d785=\"],\n  [\"66\u00d795=\", \"48\u00d730=\"],\n  [\"80\u00d746=\", \"58\u00d746=\"],\n  [\"78\u00d721=\", \"75\u00d713=\"],\n  [\"98\u00d718=\", \"83\u00d715=\"],\n  [\"76\u00d712=\", \"98\u00d787=\"],\n  [\"94\u00d780=\", \"46\u00d749=\"],\n  [\"69\u00d791=\", \"68\u00d791=\"],\n  [\"98\u00d715=\", \"94\u00d720=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and every two-digit-by-two-digit multiplication\n# problem in the practice sheet with the new values from the target\n# revision. Every \"before\" string is unique within the document, so a\n# Find/Replace pass per pair (wdReplaceAll semantics, but each needle\n# only ever matches once) is unambiguous.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{old=\"2025-08-09 Saturday\"; new=\"2025-08-10 Sunday\"},\n    @{old=\"92\u00d740=\"; new=\"61\u00d734=\"},\n    @{old=\"14\u00d773=\"; new=\"19\u00d745=\"},\n    @{old=\"34\u00d716=\"; new=\"78\u00d778=\"},\n    @{old=\"49\u00d721=\"; new=\"62\u00d713=\"},\n    @{old=\"71\u00d768=\"; new=\"92\u00d790=\"},\n    @{old=\"39\u00d781=\"; new=\"83\u00d772=\"},\n    @{old=\"74\u00d782=\"; new=\"95\u00d740=\"},\n    @{old=\"85\u00d798=\"; new=\"96\u00d721=\"},\n    @{old=\"16\u00d766=\"; new=\"40\u00d731=\"},\n    @{old=\"34\u00d781=\"; new=\"58\u00d756=\"},\n    @{old=\"79\u00d773=\"; new=\"42\u00d748=\"},\n    @{old=\"62\u00d764=\"; new=\"53\u00d732=\"},\n    @{old=\"53\u00d737=\"; new=\"81\u00d790=\"},\n    @{old=\"47\u00d745=\"; new=\"73\u00d785=\"},\n    @{old=\"62\u00d721=\"; new=\"18\u00d765=\"},\n    @{old=\"45\u00d711=\"; new=\"13\u00d771=\"},\n    @{old=\"59\u00d756=\"; new=\"49\u00d785=\"},\n    @{old=\"66\u00d795=\"; new=\"48\u00d730=\"},\n    @{old=\"80\u00d746=\"; new=\"58\u00d746=\"},\n    @{old=\"78\u00d721=\"; new=\"75\u00d713=\"},\n    @{old=\"98\u00d718=\"; new=\"83\u00d715=\"},\n    @{old=\"76\u00d712=\"; new=\"98\u00d787=\"},\n    @{old=\"94\u00d780=\"; new=\"46\u00d749=\"},\n    @{old=\"69\u00d791=\"; new=\"68\u00d791=\"},\n    @{old=\"98\u00d715=\"; new=\"94\u00d720=\"}\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.old\n    $find.Replacement.Text = $r.new\n    $find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)\n}\n"}
